$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns D need to be forced to text so Excel does not reinterpret
# numeric-looking strings (e.g. "21.50", "27.346.09", "0.00001030")
# as numbers and strip formatting / precision.
$textCells = @('D2', 'D3', 'D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D22', 'D23', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.346.09'
$ws.Range('E2').Value = '  -3.29%  '

$ws.Range('D3').Value = '1.856.82'
$ws.Range('E3').Value = '  -3.31%  '

$ws.Range('E4').Value = '  +0.38%  '

$ws.Range('D5').Value = '329.02'
$ws.Range('E5').Value = '  +0.08%  '

$ws.Range('E6').Value = '  +0.25%  '

$ws.Range('D7').Value = '0.4613'
$ws.Range('E7').Value = '  -1.74%  '

$ws.Range('D8').Value = '0.3949'
$ws.Range('E8').Value = '  -2.06%  '

$ws.Range('D9').Value = '46.46'
$ws.Range('E9').Value = '  -12.54%  '

$ws.Range('D10').Value = '0.07936'
$ws.Range('E10').Value = '  -5.93%  '

$ws.Range('D11').Value = '1.012'
$ws.Range('E11').Value = '  -3.55%  '

$ws.Range('D12').Value = '21.50'
$ws.Range('E12').Value = '  -3.45%  '

$ws.Range('D13').Value = '1.845.78'
$ws.Range('E13').Value = '  -3.81%  '

$ws.Range('D14').Value = '5.928'
$ws.Range('E14').Value = '  -2.61%  '

$ws.Range('D15').Value = '7.144'
$ws.Range('E15').Value = '  -4.02%  '

$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  +0.40%  '

$ws.Range('D17').Value = '86.20'
$ws.Range('E17').Value = '  -4.14%  '

$ws.Range('D18').Value = '0.00001030'
$ws.Range('E18').Value = '  -3.44%  '

$ws.Range('D19').Value = '0.06583'
$ws.Range('E19').Value = '  -0.43%  '

$ws.Range('D20').Value = '17.25'
$ws.Range('E20').Value = '  -4.60%  '

$ws.Range('E21').Value = '  +0.18%  '

$ws.Range('D22').Value = '5.481'
$ws.Range('E22').Value = '  -4.37%  '

$ws.Range('D23').Value = '27.354.99'
$ws.Range('E23').Value = '  -3.23%  '

$ws.Range('E24').Value = '  -3.50%  '

$ws.Range('E25').Value = '  +1.01%  '

$ws.Range('D26').Value = '2.071.91'
$ws.Range('E26').Value = '  -3.48%  '

$ws.Range('D27').Value = '153.46'
$ws.Range('E27').Value = '  +0.04%  '

$ws.Range('D28').Value = '20.21'
$ws.Range('E28').Value = '  +0.77%  '

$ws.Range('D29').Value = '2.064'
$ws.Range('E29').Value = '  -3.05%  '

$ws.Range('D30').Value = '5.470'
$ws.Range('E30').Value = '  -5.03%  '

$ws.Range('D31').Value = '121.67'
$ws.Range('E31').Value = '  -1.65%  '

$ws.Range('D32').Value = '0.09421'
$ws.Range('E32').Value = '  -2.22%  '

$ws.Range('D33').Value = '0.9492'
$ws.Range('E33').Value = '  -3.05%  '

$ws.Range('D34').Value = '1.446'
$ws.Range('E34').Value = '  -1.06%  '

$ws.Range('D35').Value = '3.590'
$ws.Range('E35').Value = '  -1.19%  '

$ws.Range('D36').Value = '5.265'
$ws.Range('E36').Value = '  -5.26%  '

$ws.Range('D37').Value = '0.06041'
$ws.Range('E37').Value = '  -2.35%  '

$ws.Range('D38').Value = '0.02227'
$ws.Range('E38').Value = '  -3.42%  '

$ws.Range('D39').Value = '1.213'
$ws.Range('E39').Value = '  -3.89%  '

$ws.Range('B40').Value = 'Frax'
$ws.Range('C40').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D40').Value = '1.002'
$ws.Range('E40').Value = '  +0.15%  '

$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '8.035'
$ws.Range('E41').Value = '  -8.88%  '

$ws.Range('D42').Value = '0.5924'
$ws.Range('E42').Value = '  -3.77%  '

$ws.Range('D43').Value = '0.1885'
$ws.Range('E43').Value = '  -1.25%  '

$ws.Range('D44').Value = '10.18'
$ws.Range('E44').Value = '  -8.20%  '

$ws.Range('D45').Value = '1.285'
$ws.Range('E45').Value = '  -1.46%  '

$ws.Range('D46').Value = '0.5626'
$ws.Range('E46').Value = '  -4.08%  '

$ws.Range('D47').Value = '12.01'
$ws.Range('E47').Value = '  -6.24%  '

$ws.Range('D48').Value = '3.397'
$ws.Range('E48').Value = '  -1.14%  '

$ws.Range('D49').Value = '1.916'
$ws.Range('E49').Value = '  -5.73%  '

$ws.Range('D50').Value = '0.06764'
$ws.Range('E50').Value = '  -2.02%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.00000000309'
$ws.Range('E51').Value = '  +2.49%  '

# Restore default (General) style on the cells we forced to text,
# so the resulting cells carry no explicit style index, matching
# the original workbook formatting.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
